$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Generate Report for Handoff" — a new handoff (xliff generation) round ran
# for the two files that previously showed "Handed back: in sync with en-US"
# (a85fc7a7-... and e7c7e19a-...). Their status flips to "Ready for handoff",
# a fresh handoff timestamp is recorded, and — because the handback on file
# is now stale relative to the newest source commit — an Error Detail note
# is attached on each locale sheet.
# ---------------------------------------------------------------------------

$newStatus = "Ready for handoff"

# ---------------------------------------------------------------------------
# Overview sheet: rows for a85fc7a7 (row 4) and e7c7e19a (row 5)
#   E/F (zh-cn / de-de status columns) -> "Ready for handoff"
#   G (Latest HO Xliff Generate Date)  -> "2016-08-30 22:28:23"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E4").Value = $newStatus
$wsOverview.Range("F4").Value = $newStatus
$wsOverview.Range("G4").Value = "2016-08-30 22:28:23"

$wsOverview.Range("E5").Value = $newStatus
$wsOverview.Range("F5").Value = $newStatus
$wsOverview.Range("G5").Value = "2016-08-30 22:28:23"

# ---------------------------------------------------------------------------
# zh-cn sheet: rows 4 (a85fc7a7) and 5 (e7c7e19a)
#   C  (Status)                  -> "Ready for handoff"
#   H  (Latest Handoff Datetime) -> "2016-08-30 22:28:18"
#   P  (Error Detail)            -> stale-handback warning
#   Column P widened to fit the long message
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C4").Value = $newStatus
$wsZhCn.Range("H4").Value = "2016-08-30 22:28:18"
$wsZhCn.Range("P4").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0e74b191ea76327f91b0dc59b11b209345ef2053/e2e/a85fc7a7-8339-4231-8b93-afe5b823b301.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9559f97ff718bfebed012ad11a38aabbae5284c4/e2e/a85fc7a7-8339-4231-8b93-afe5b823b301.md."

$wsZhCn.Range("C5").Value = $newStatus
$wsZhCn.Range("H5").Value = "2016-08-30 22:28:18"
$wsZhCn.Range("P5").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0e74b191ea76327f91b0dc59b11b209345ef2053/e2e/e7c7e19a-de4c-446d-8d48-bc073ba7b840.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9559f97ff718bfebed012ad11a38aabbae5284c4/e2e/e7c7e19a-de4c-446d-8d48-bc073ba7b840.md."

$wsZhCn.Columns.Item(16).ColumnWidth = 39.1640625

# ---------------------------------------------------------------------------
# de-de sheet: rows 4 (a85fc7a7) and 5 (e7c7e19a)
#   C  (Status)                  -> "Ready for handoff"
#   H  (Latest Handoff Datetime) -> "2016-08-30 22:28:23"
#   P  (Error Detail)            -> stale-handback warning
#   Column P widened to fit the long message
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C4").Value = $newStatus
$wsDeDe.Range("H4").Value = "2016-08-30 22:28:23"
$wsDeDe.Range("P4").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0e74b191ea76327f91b0dc59b11b209345ef2053/e2e/a85fc7a7-8339-4231-8b93-afe5b823b301.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9559f97ff718bfebed012ad11a38aabbae5284c4/e2e/a85fc7a7-8339-4231-8b93-afe5b823b301.md."

$wsDeDe.Range("C5").Value = $newStatus
$wsDeDe.Range("H5").Value = "2016-08-30 22:28:23"
$wsDeDe.Range("P5").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0e74b191ea76327f91b0dc59b11b209345ef2053/e2e/e7c7e19a-de4c-446d-8d48-bc073ba7b840.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9559f97ff718bfebed012ad11a38aabbae5284c4/e2e/e7c7e19a-de4c-446d-8d48-bc073ba7b840.md."

$wsDeDe.Columns.Item(16).ColumnWidth = 39.1640625
